# Updated cryptos list on Sun Oct  1 18:27:12 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) / Volume(1h) (column E) snapshot values for
# every coin row on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price strings parse as a clean decimal number
# (e.g. "214.20", "0.518"). Those cells are pre-marked as Text so the
# COM layer stores the literal digits/trailing zeros instead of silently
# coercing them to a Double (which would turn "214.20" into 214.2).
$textFormatCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D30",
    "D34",
    "D36",
    "D37",
    "D39",
    "D41",
    "D42",
    "D46",
    "D47",
    "D50"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.147.94"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.678.83"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "214.20"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "22.74"
$ws.Range("E8").Value = "  +6.30%  "
$ws.Range("D9").Value = "0.260"
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").Value = "0.0622"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "1.915.34"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.664.63"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "4.19"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").Value = "0.550"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "66.56"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "27.102.38"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "235.62"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "7.89"
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("D20").Value = "0.0₃0740"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "4.53"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "146.85"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "7.40"
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "1.541.81"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "3.24"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").Value = "0.605"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").Value = "0.940"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "0.0172"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "69.27"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.822.67"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "0.777"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "89.82"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("D50").Value = "8.23"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("E51").Value = "  +0.06%  "
